$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Rows 110-112 (match rows for 109, 110, 111 in the "Indice" column) get
# cyclically rotated: the row that used to sit at 110 moves to 112, the row
# that used to sit at 111 moves to 110, and the row that used to sit at 112
# moves to 111. The "Indice"(A), pais(B), torneio(C), temporada(D),
# data_partida(E) and the four opening-odds timestamp columns (K, O, S) stay
# put because they're identical across the three rows; everything else
# (teams, scores, odds, closing timestamps, url) travels with its match.
# ---------------------------------------------------------------------------

function Get-RowData($row) {
    $cols = @("F","G","H","I","J","L","M","N","P","Q","R","T","U","V")
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range($col + $row).Value2
    }
    return $data
}

$row110 = Get-RowData 110
$row111 = Get-RowData 111
$row112 = Get-RowData 112

function Set-RowData($row, $data) {
    foreach ($col in $data.Keys) {
        $ws.Range($col + $row).Value = $data[$col]
    }
}

# new row 110 <- old row 111, new row 111 <- old row 112, new row 112 <- old row 110
Set-RowData 110 $row111
Set-RowData 111 $row112
Set-RowData 112 $row110

# ---------------------------------------------------------------------------
# Two brand-new match rows appended at the bottom (114, 115). Copy the
# formatting from row 113 (the last existing data row) first so the new
# rows pick up the same styles (bold/boxed index column, date-time number
# format on data_partida), then fill in the values.
# ---------------------------------------------------------------------------

$ws.Range("A113:V113").Copy() | Out-Null
$ws.Range("A114:V115").PasteSpecial(-4122) | Out-Null

$ws.Range("A114").Value = 113
$ws.Range("B114").Value = "ecuador"
$ws.Range("C114").Value = "liga-pro"
$ws.Range("D114").Value = "'2023"
$ws.Range("E114").Value = 45270.9375
$ws.Range("F114").Value = "Ind. del Valle"
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = "LDU Quito"
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 2.06
$ws.Range("K114").Value = "09/12/2023 11:43"
$ws.Range("L114").Value = 2.33
$ws.Range("M114").Value = "10/12/2023 22:20"
$ws.Range("N114").Value = 3.4
$ws.Range("O114").Value = "09/12/2023 11:43"
$ws.Range("P114").Value = 3.1
$ws.Range("Q114").Value = "10/12/2023 22:20"
$ws.Range("R114").Value = 3.37
$ws.Range("S114").Value = "09/12/2023 11:43"
$ws.Range("T114").Value = 3.42
$ws.Range("U114").Value = "10/12/2023 22:20"
$ws.Range("V114").Value = "https://www.betexplorer.com/football/ecuador/liga-pro/independiente-del-valle-ldu-quito/joyDcvi4/"

$ws.Range("A115").Value = 114
$ws.Range("B115").Value = "ecuador"
$ws.Range("C115").Value = "liga-pro"
$ws.Range("D115").Value = "'2023"
$ws.Range("E115").Value = 45277.9375
$ws.Range("F115").Value = "LDU Quito"
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = "Ind. del Valle"
$ws.Range("I115").Value = 1
$ws.Range("J115").Value = 2.18
$ws.Range("K115").Value = "10/12/2023 22:42"
$ws.Range("L115").Value = 2.66
$ws.Range("M115").Value = "17/12/2023 22:27"
$ws.Range("N115").Value = 3.36
$ws.Range("O115").Value = "10/12/2023 22:42"
$ws.Range("P115").Value = 2.98
$ws.Range("Q115").Value = "17/12/2023 22:24"
$ws.Range("R115").Value = 3.14
$ws.Range("S115").Value = "10/12/2023 22:42"
$ws.Range("T115").Value = 3.01
$ws.Range("U115").Value = "17/12/2023 22:27"
$ws.Range("V115").Value = "https://www.betexplorer.com/football/ecuador/liga-pro/ldu-quito-independiente-del-valle/p8zHdb6A/"
